$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = "AS A Buyer I Should be able to Register to login further So That it helps me in historical comparision."
$ws.Range("C25").Interior.Color = 65535
$ws.Range("A26").Select()
